$wb = $excel.ActiveWorkbook

# --- Rename existing "Robot Info" sheet to "Robot Info-old" ---
$robotOld = $wb.Worksheets.Item("Robot Info")
$robotOld.Name = "Robot Info-old"

# --- Add new "Robot Info" sheet right after "Robot Info-old" ---
$robotNew = $wb.Worksheets.Add($null, $robotOld)
$robotNew.Name = "Robot Info"

# Fill in headers A1:C1 first (matches authoring order so shared strings line up)
$robotNew.Range("A1").Value = "Robot"
$robotNew.Range("B1").Value = "Type"
$robotNew.Range("C1").Value = "Cost"

# --- Add new "Package Configs" sheet right after the new "Robot Info" ---
$pkg = $wb.Worksheets.Add($null, $robotNew)
$pkg.Name = "Package Configs"

$pkg.Range("B1").Value = "Package 1"
$pkg.Range("C1").Value = "Package 2"
$pkg.Range("D1").Value = "Package 3"
$pkg.Range("E1").Value = "Package 4"
$pkg.Range("F1").Value = "Package 5"
$pkg.Range("G1").Value = "Package 6"

$pkg.Range("A2").Value = "Robot 1"
$pkg.Range("A3").Value = "Robot 2"
$pkg.Range("A4").Value = "Robot 3"
$pkg.Range("A5").Value = "Robot 4"
$pkg.Range("A6").Value = "Robot 5"
$pkg.Range("A7").Value = "Profit"

$pkgData = @(
    @(1,2,1,2,0,1),
    @(0,1,0,2,2,2),
    @(2,2,1,3,2,3),
    @(2,2,2,3,3,4),
    @(4,4,4,5,3,3),
    @(20,17,18,16,19,20)
)
$r = 2
foreach ($row in $pkgData) {
    $pkg.Cells.Item($r, 2).Value = $row[0]
    $pkg.Cells.Item($r, 3).Value = $row[1]
    $pkg.Cells.Item($r, 4).Value = $row[2]
    $pkg.Cells.Item($r, 5).Value = $row[3]
    $pkg.Cells.Item($r, 6).Value = $row[4]
    $pkg.Cells.Item($r, 7).Value = $row[5]
    $r++
}

# --- Back to "Robot Info": add the "Available" column + row data ---
$robotNew.Range("D1").Value = "Available"

$robotData = @(
    @(1, "Vaccum", 400, 35),
    @(2, "Carpet Cleaner", 600, 25),
    @(3, "Sink Cleaner", 500, 20),
    @(4, "Duster Drone", 300, 30),
    @(5, "Window Cleaner", 600, 35)
)
$r = 2
foreach ($row in $robotData) {
    $robotNew.Cells.Item($r, 1).Value = $row[0]
    $robotNew.Cells.Item($r, 2).Value = $row[1]
    $robotNew.Cells.Item($r, 3).Value = $row[2]
    $robotNew.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Apply the built-in "Currency" cell style to the Cost column (C1:C6)
$robotNew.Range("C1:C6").Style = "Currency"

# --- Column widths (best-fit, as Excel would auto-size after typing) ---
$robotNew.Columns.Item(2).AutoFit()
$robotNew.Columns.Item(3).AutoFit()
$pkg.Columns.Item(3).AutoFit()

# --- Selections left by the author in each sheet ---
[void]$robotNew.Range("D7").Select()
[void]$pkg.Range("G8").Select()

# Package Configs ends up the active/visible tab
[void]$pkg.Activate()
